$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.137.83'
$ws.Range("E2").Value = '  +9.10%  '
$ws.Range("D3").Value = '3.222.62'
$ws.Range("E3").Value = '  +3.77%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.Value = "'397.29"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.25%  '
$c = $ws.Range("D6")
$c.Value = "'110.71"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +6.41%  '
$c = $ws.Range("D7")
$c.Value = "'0.552"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +2.49%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +4.91%  '
$c = $ws.Range("D10")
$c.Value = "'39.25"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +5.48%  '
$ws.Range("E11").Value = '  +6.12%  '
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("D13").Value = '3.740.01'
$ws.Range("E13").Value = '  +3.97%  '
$c = $ws.Range("D14")
$c.Value = "'8.08"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +3.63%  '
$c = $ws.Range("D15")
$c.Value = "'19.03"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.37%  '
$ws.Range("D16").Value = '3.218.91'
$ws.Range("E16").Value = '  +3.63%  '
$ws.Range("E17").Value = '  +4.95%  '
$c = $ws.Range("D18")
$c.Value = "'10.84"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("D19").Value = '56.038.19'
$ws.Range("E20").Value = '  +1.80%  '
$ws.Range("E21").Value = '  +5.91%  '
$c = $ws.Range("D22")
$c.Value = "'13.04"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.99%  '
$c = $ws.Range("D23")
$c.Value = "'299.49"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +12.20%  '
$c = $ws.Range("D24")
$c.Value = "'75.44"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +7.25%  '
$c = $ws.Range("D25")
$c.Value = "'3.22"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.34%  '
$c = $ws.Range("D27")
$c.Value = "'28.16"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.57%  '
$c = $ws.Range("D28")
$c.Value = "'7.47"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("E29").Value = '  +3.95%  '
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("E31").Value = '  +3.30%  '
$c = $ws.Range("D32")
$c.Value = "'11.15"
$c.Style = "Normal"
$ws.Range("E32").Value = '  +6.74%  '
$c = $ws.Range("D33")
$c.Value = "'0.0491"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +3.37%  '
$c = $ws.Range("D34")
$c.Value = "'36.14"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.82%  '
$c = $ws.Range("D35")
$c.Value = "'2.13"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +2.37%  '
$c = $ws.Range("D36")
$c.Value = "'51.35"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.55%  '
$ws.Range("E37").Value = '  +25.38%  '
$ws.Range("E38").Value = '  +4.00%  '
$c = $ws.Range("D39")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.08%  '
$c = $ws.Range("D40")
$c.Value = "'134.54"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.90%  '
$ws.Range("E41").Value = '  +3.11%  '
$c = $ws.Range("D42")
$c.Value = "'17.36"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +4.57%  '
$ws.Range("E43").Value = '  +4.88%  '
$ws.Range("E44").Value = '  +2.85%  '
$c = $ws.Range("D45")
$c.Value = "'0.282"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.58%  '
$c = $ws.Range("D46")
$c.Value = "'22.23"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.01%  '
$c = $ws.Range("D47")
$c.Value = "'2.15"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +49.13%  '
$c = $ws.Range("D49")
$c.Value = "'2.46"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.01%  '
$ws.Range("D50").Value = '2.132.25'
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("E51").Value = '  +9.58%  '
